# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column, inserted
# right after the existing "total" column and before the "date" column.
# Every data row gets the literal value "stock" in that new column.
#
# Old layout: name | owner | quantity | face_value | currency | total | date | legislator_name | legislator_id
# New layout: name | owner | quantity | face_value | currency | total | property_category | date | legislator_name | legislator_id

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# "total" is column G (7th column), "date" was column H (8th column).
# Insert a whole column at H so everything from H onward (date,
# legislator_name, legislator_id) shifts one column right to I/J/K,
# and the newly inserted column H inherits formatting (header bold/border
# style, data-row style) from the columns around it.
$ws.Columns.Item(8).Insert()

# Header for the freshly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Fill every data row (rows 2-11) with the literal category value.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 11 }
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
